$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 177
$ws.Range("E177").Value = "Giannina"
$ws.Range("F177").Value = "Atromitos Athinon"
$ws.Range("K177").Value = "D"
$ws.Range("B177").Value = 6937271
$ws.Range("G177").Value = 1
$ws.Range("H177").Value = 1
$ws.Range("I177").Value = 1
$ws.Range("J177").Value = 0
$ws.Range("L177").Value = 2.45
$ws.Range("M177").Value = 3.1
$ws.Range("N177").Value = 3.1
$ws.Range("O177").Value = 2
$ws.Range("P177").Value = 3.3
$ws.Range("Q177").Value = 4
$ws.Range("R177").Value = -0.5
$ws.Range("S177").Value = 2.025
$ws.Range("T177").Value = 1.825
$ws.Range("U177").Value = 2.25
$ws.Range("V177").Value = 1.85
$ws.Range("W177").Value = 2
$ws.Range("X177").Value = -1
$ws.Range("Y177").Value = 2.3
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = -1
$ws.Range("AB177").Value = 0.825
$ws.Range("AC177").Value = -0.5
$ws.Range("AD177").Value = 0.5

# Row 178
$ws.Range("E178").Value = "Kifisias FC"
$ws.Range("F178").Value = "Panetolikos"
$ws.Range("K178").Value = "D"
$ws.Range("B178").Value = 6935701
$ws.Range("G178").Value = 2
$ws.Range("H178").Value = 2
$ws.Range("I178").Value = 1
$ws.Range("J178").Value = 0
$ws.Range("L178").Value = 2.45
$ws.Range("M178").Value = 3.25
$ws.Range("N178").Value = 3
$ws.Range("O178").Value = 2.05
$ws.Range("P178").Value = 3.3
$ws.Range("Q178").Value = 3.8
$ws.Range("R178").Value = -0.5
$ws.Range("S178").Value = 2.05
$ws.Range("T178").Value = 1.8
$ws.Range("U178").Value = 2.25
$ws.Range("V178").Value = 1.8
$ws.Range("W178").Value = 2.05
$ws.Range("X178").Value = -1
$ws.Range("Y178").Value = 2.3
$ws.Range("Z178").Value = -1
$ws.Range("AA178").Value = -1
$ws.Range("AB178").Value = 0.8
$ws.Range("AC178").Value = 0.8
$ws.Range("AD178").Value = -1

# Row 179
$ws.Range("E179").Value = "Lamia"
$ws.Range("F179").Value = "PAOK Salonika"
$ws.Range("K179").Value = "A"
$ws.Range("B179").Value = 6937272
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 2
$ws.Range("I179").Value = 0
$ws.Range("J179").Value = 0
$ws.Range("L179").Value = 7.5
$ws.Range("M179").Value = 4.5
$ws.Range("N179").Value = 1.444
$ws.Range("O179").Value = 9.5
$ws.Range("P179").Value = 5
$ws.Range("Q179").Value = 1.333
$ws.Range("R179").Value = 1.5
$ws.Range("S179").Value = 1.925
$ws.Range("T179").Value = 1.925
$ws.Range("U179").Value = 3
$ws.Range("V179").Value = 1.95
$ws.Range("W179").Value = 1.9
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = -1
$ws.Range("Z179").Value = 0.333
$ws.Range("AA179").Value = -1
$ws.Range("AB179").Value = 0.925
$ws.Range("AC179").Value = -1
$ws.Range("AD179").Value = 0.8999999999999999

# Row 180
$ws.Range("E180").Value = "Panserraikos"
$ws.Range("F180").Value = "Asteras Tripolis"
$ws.Range("K180").Value = "H"
$ws.Range("B180").Value = 6935700
$ws.Range("G180").Value = 2
$ws.Range("H180").Value = 1
$ws.Range("I180").Value = 1
$ws.Range("J180").Value = 1
$ws.Range("L180").Value = 2.6
$ws.Range("M180").Value = 3.2
$ws.Range("N180").Value = 2.875
$ws.Range("O180").Value = 2.25
$ws.Range("P180").Value = 3.3
$ws.Range("Q180").Value = 3.3
$ws.Range("R180").Value = -0.25
$ws.Range("S180").Value = 1.925
$ws.Range("T180").Value = 1.925
$ws.Range("U180").Value = 2.25
$ws.Range("V180").Value = 2
$ws.Range("W180").Value = 1.85
$ws.Range("X180").Value = 1.25
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = -1
$ws.Range("AA180").Value = 0.925
$ws.Range("AB180").Value = -1
$ws.Range("AC180").Value = 1
$ws.Range("AD180").Value = -1

# Row 181
$ws.Range("E181").Value = "Olympiakos"
$ws.Range("F181").Value = "Volos NFC"
$ws.Range("K181").Value = "H"
$ws.Range("B181").Value = 6937270
$ws.Range("G181").Value = 3
$ws.Range("H181").Value = 0
$ws.Range("I181").Value = 2
$ws.Range("J181").Value = 0
$ws.Range("L181").Value = 1.125
$ws.Range("M181").Value = 9
$ws.Range("N181").Value = 19
$ws.Range("O181").Value = 1.111
$ws.Range("P181").Value = 9
$ws.Range("Q181").Value = 21
$ws.Range("R181").Value = -2.25
$ws.Range("S181").Value = 1.875
$ws.Range("T181").Value = 1.975
$ws.Range("U181").Value = 3.25
$ws.Range("V181").Value = 2
$ws.Range("W181").Value = 1.85
$ws.Range("X181").Value = 0.111
$ws.Range("Y181").Value = -1
$ws.Range("Z181").Value = -1
$ws.Range("AA181").Value = 0.875
$ws.Range("AB181").Value = -1
$ws.Range("AC181").Value = -0.5
$ws.Range("AD181").Value = 0.425

# Row 222
$ws.Range("E222").Value = "Kifisias FC"
$ws.Range("F222").Value = "Giannina"
$ws.Range("K222").Value = "A"
$ws.Range("B222").Value = 7920464
$ws.Range("G222").Value = 2
$ws.Range("H222").Value = 3
$ws.Range("I222").Value = 1
$ws.Range("J222").Value = 1
$ws.Range("L222").Value = 1.571
$ws.Range("M222").Value = 3.8
$ws.Range("N222").Value = 6.5
$ws.Range("O222").Value = 1.4
$ws.Range("P222").Value = 4.75
$ws.Range("Q222").Value = 7.5
$ws.Range("R222").Value = -1.25
$ws.Range("S222").Value = 2.05
$ws.Range("T222").Value = 1.8
$ws.Range("U222").Value = 2.75
$ws.Range("V222").Value = 1.925
$ws.Range("W222").Value = 1.925
$ws.Range("X222").Value = -1
$ws.Range("Y222").Value = -1
$ws.Range("Z222").Value = 6.5
$ws.Range("AA222").Value = -1
$ws.Range("AB222").Value = 0.8
$ws.Range("AC222").Value = 0.925
$ws.Range("AD222").Value = -1

# Row 223
$ws.Range("E223").Value = "Asteras Tripolis"
$ws.Range("F223").Value = "OFI Crete"
$ws.Range("K223").Value = "D"
$ws.Range("B223").Value = 7920462
$ws.Range("G223").Value = 1
$ws.Range("H223").Value = 1
$ws.Range("I223").Value = 1
$ws.Range("J223").Value = 0
$ws.Range("L223").Value = 2.5
$ws.Range("M223").Value = 3
$ws.Range("N223").Value = 3.1
$ws.Range("O223").Value = 2.55
$ws.Range("P223").Value = 2.75
$ws.Range("Q223").Value = 3.3
$ws.Range("R223").Value = -0.25
$ws.Range("S223").Value = 2.125
$ws.Range("T223").Value = 1.75
$ws.Range("U223").Value = 2.25
$ws.Range("V223").Value = 1.85
$ws.Range("W223").Value = 2
$ws.Range("X223").Value = -1
$ws.Range("Y223").Value = 1.75
$ws.Range("Z223").Value = -1
$ws.Range("AA223").Value = -0.5
$ws.Range("AB223").Value = 0.375
$ws.Range("AC223").Value = -0.5
$ws.Range("AD223").Value = 0.5
